# Auto-generated update of market-price derived columns (H:N) across multiple
# profession sheets, refreshed by the scheduled Jenova_Profits runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Cells.Item(11, 8).Value = 196  # H11: 158.25 -> 196
$ws.Cells.Item(11, 9).Value = 196  # I11: 158.25 -> 196
$ws.Cells.Item(11, 11).Value = 196  # K11: 158.25 -> 196
$ws.Cells.Item(11, 13).Value = -56  # M11: -18.25 -> -56

# Row 18
$ws.Cells.Item(18, 8).Value = 320.1  # H18: 335.66666 -> 320.1
$ws.Cells.Item(18, 9).Value = 300.1111  # I18: 315.125 -> 300.1111
$ws.Cells.Item(18, 11).Value = 300.1111  # K18: 315.125 -> 300.1111
$ws.Cells.Item(18, 13).Value = -16.11110000000002  # M18: -31.125 -> -16.11110000000002

# Row 40
$ws.Cells.Item(40, 8).Value = 8250.071  # H40: 8772.817999999999 -> 8250.071
$ws.Cells.Item(40, 10).Value = 10570.714  # J40: 13748.75 -> 10570.714
$ws.Cells.Item(40, 12).Value = 10570.714  # L40: 13748.75 -> 10570.714
$ws.Cells.Item(40, 14).Value = -10920.714  # N40: -14098.75 -> -10920.714

# Row 43
$ws.Cells.Item(43, 8).Value = 1534.8235  # H43: 2086.3635 -> 1534.8235
$ws.Cells.Item(43, 9).Value = 1538.4615  # I43: 1883.3334 -> 1538.4615
$ws.Cells.Item(43, 10).Value = 1523  # J43: 3000 -> 1523
$ws.Cells.Item(43, 11).Value = 1538.4615  # K43: 1883.3334 -> 1538.4615
$ws.Cells.Item(43, 12).Value = 1523  # L43: 3000 -> 1523
$ws.Cells.Item(43, 13).Value = -1469.4615  # M43: -1814.3334 -> -1469.4615
$ws.Cells.Item(43, 14).Value = -1661  # N43: -3138 -> -1661

# Row 88
$ws.Cells.Item(88, 8).Value = 2010.3  # H88: 2022.6666 -> 2010.3
$ws.Cells.Item(88, 10).Value = 2338.125  # J88: 2400.8572 -> 2338.125
$ws.Cells.Item(88, 12).Value = 2338.125  # L88: 2400.8572 -> 2338.125
$ws.Cells.Item(88, 14).Value = -3150.125  # N88: -3212.8572 -> -3150.125

# Row 91
$ws.Cells.Item(91, 8).Value = 2010.3  # H91: 2022.6666 -> 2010.3
$ws.Cells.Item(91, 10).Value = 2338.125  # J91: 2400.8572 -> 2338.125
$ws.Cells.Item(91, 12).Value = 2338.125  # L91: 2400.8572 -> 2338.125
$ws.Cells.Item(91, 14).Value = -5146.125  # N91: -5208.8572 -> -5146.125

# Row 94
$ws.Cells.Item(94, 8).Value = 1910.5834  # H94: 2075.2727 -> 1910.5834
$ws.Cells.Item(94, 9).Value = 1420.1111  # I94: 1572.75 -> 1420.1111
$ws.Cells.Item(94, 10).Value = 3382  # J94: 3415.3333 -> 3382
$ws.Cells.Item(94, 11).Value = 1420.1111  # K94: 1572.75 -> 1420.1111
$ws.Cells.Item(94, 12).Value = 3382  # L94: 3415.3333 -> 3382
$ws.Cells.Item(94, 13).Value = -969.1111000000001  # M94: -1121.75 -> -969.1111000000001
$ws.Cells.Item(94, 14).Value = -4284  # N94: -4317.3333 -> -4284

# Row 116
$ws.Cells.Item(116, 8).Value = 17668.666  # H116: 19139.375 -> 17668.666
$ws.Cells.Item(116, 9).Value = 7123  # I116: 7422 -> 7123
$ws.Cells.Item(116, 10).Value = 30850.75  # J116: 38668.332 -> 30850.75
$ws.Cells.Item(116, 11).Value = 7123  # K116: 7422 -> 7123
$ws.Cells.Item(116, 12).Value = 30850.75  # L116: 38668.332 -> 30850.75
$ws.Cells.Item(116, 13).Value = -3681  # M116: -3980 -> -3681
$ws.Cells.Item(116, 14).Value = -37734.75  # N116: -45552.332 -> -37734.75

# Row 132
$ws.Cells.Item(132, 8).Value = 7457.5  # H132: 7569.154 -> 7457.5
$ws.Cells.Item(132, 9).Value = 1357.772  # I132: 1378.4464 -> 1357.772
$ws.Cells.Item(132, 11).Value = 4073.316  # K132: 4135.3392 -> 4073.316
$ws.Cells.Item(132, 13).Value = -1543.316  # M132: -1605.3392 -> -1543.316

# Row 133
$ws.Cells.Item(133, 8).Value = 38840.8  # H133: 38582.547 -> 38840.8
$ws.Cells.Item(133, 10).Value = 38840.8  # J133: 38582.547 -> 38840.8
$ws.Cells.Item(133, 12).Value = 38840.8  # L133: 38582.547 -> 38840.8
$ws.Cells.Item(133, 14).Value = -48960.8  # N133: -48702.547 -> -48960.8

# Row 134
$ws.Cells.Item(134, 8).Value = 49822  # H134: 89033 -> 49822
$ws.Cells.Item(134, 10).Value = 49822  # J134: 89033 -> 49822
$ws.Cells.Item(134, 12).Value = 49822  # L134: 89033 -> 49822
$ws.Cells.Item(134, 14).Value = -59962  # N134: -99173 -> -59962

# Row 135
$ws.Cells.Item(135, 8).Value = 627110.9  # H135: 590537.9 -> 627110.9
$ws.Cells.Item(135, 9).Value = 716153.0600000001  # I135: 771116.4 -> 716153.0600000001
$ws.Cells.Item(135, 10).Value = 3815.5  # J135: 3657.75 -> 3815.5
$ws.Cells.Item(135, 11).Value = 6445377.540000001  # K135: 6940047.600000001 -> 6445377.540000001
$ws.Cells.Item(135, 12).Value = 34339.5  # L135: 32919.75 -> 34339.5
$ws.Cells.Item(135, 13).Value = -6442842.540000001  # M135: -6937512.600000001 -> -6442842.540000001
$ws.Cells.Item(135, 14).Value = -39409.5  # N135: -37989.75 -> -39409.5

# Row 137
$ws.Cells.Item(137, 8).Value = 4501.1  # H137: 4949.25 -> 4501.1
$ws.Cells.Item(137, 9).Value = 4851.2085  # I137: 5165.625 -> 4851.2085
$ws.Cells.Item(137, 10).Value = 3100.6667  # J137: 3651 -> 3100.6667
$ws.Cells.Item(137, 11).Value = 14553.6255  # K137: 15496.875 -> 14553.6255
$ws.Cells.Item(137, 12).Value = 9302.000100000001  # L137: 10953 -> 9302.000100000001
$ws.Cells.Item(137, 13).Value = -12003.6255  # M137: -12946.875 -> -12003.6255
$ws.Cells.Item(137, 14).Value = -14402.0001  # N137: -16053 -> -14402.0001

# Row 141
$ws.Cells.Item(141, 8).Value = 3699.1428  # H141: 3829.923 -> 3699.1428
$ws.Cells.Item(141, 9).Value = 3214.8333  # I141: 3325.3635 -> 3214.8333
$ws.Cells.Item(141, 11).Value = 9644.499899999999  # K141: 9976.0905 -> 9644.499899999999
$ws.Cells.Item(141, 13).Value = -4464.499899999999  # M141: -4796.0905 -> -4464.499899999999

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Cells.Item(45, 8).Value = 2210.2  # H45: 2191.5625 -> 2210.2
$ws.Cells.Item(45, 9).Value = 1750.4546  # I45: 1763.9166 -> 1750.4546
$ws.Cells.Item(45, 11).Value = 1750.4546  # K45: 1763.9166 -> 1750.4546
$ws.Cells.Item(45, 13).Value = -1373.4546  # M45: -1386.9166 -> -1373.4546

# Row 61
$ws.Cells.Item(61, 8).Value = 3084.5862  # H61: 3014.7666 -> 3084.5862
$ws.Cells.Item(61, 9).Value = 2601.12  # I61: 2539.1538 -> 2601.12
$ws.Cells.Item(61, 11).Value = 2601.12  # K61: 2539.1538 -> 2601.12
$ws.Cells.Item(61, 13).Value = -2389.12  # M61: -2327.1538 -> -2389.12

# Row 122
$ws.Cells.Item(122, 8).Value = 7878.9473  # H122: 8317.647000000001 -> 7878.9473
$ws.Cells.Item(122, 9).Value = 19266.666  # I122: 26950 -> 19266.666
$ws.Cells.Item(122, 10).Value = 5743.75  # J122: 5833.3335 -> 5743.75
$ws.Cells.Item(122, 11).Value = 57799.99800000001  # K122: 80850 -> 57799.99800000001
$ws.Cells.Item(122, 12).Value = 17231.25  # L122: 17500.0005 -> 17231.25
$ws.Cells.Item(122, 13).Value = -55349.99800000001  # M122: -78400 -> -55349.99800000001
$ws.Cells.Item(122, 14).Value = -22131.25  # N122: -22400.0005 -> -22131.25

# Row 132
$ws.Cells.Item(132, 8).Value = 4475.75  # H132: 4535.2905 -> 4475.75
$ws.Cells.Item(132, 9).Value = 4471.7617  # I132: 4563.85 -> 4471.7617
$ws.Cells.Item(132, 11).Value = 13415.2851  # K132: 13691.55 -> 13415.2851
$ws.Cells.Item(132, 13).Value = -10885.2851  # M132: -11161.55 -> -10885.2851

# Row 136
$ws.Cells.Item(136, 8).Value = 3084.5862  # H136: 3014.7666 -> 3084.5862
$ws.Cells.Item(136, 9).Value = 2601.12  # I136: 2539.1538 -> 2601.12
$ws.Cells.Item(136, 11).Value = 7803.36  # K136: 7617.4614 -> 7803.36
$ws.Cells.Item(136, 13).Value = -5253.36  # M136: -5067.4614 -> -5253.36

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 51308.43  # H31: 51327.477 -> 51308.43
$ws.Cells.Item(31, 9).Value = 1183.8  # I31: 1263.8 -> 1183.8
$ws.Cells.Item(31, 11).Value = 1183.8  # K31: 1263.8 -> 1183.8
$ws.Cells.Item(31, 13).Value = -888.8  # M31: -968.8 -> -888.8

# Row 34
$ws.Cells.Item(34, 8).Value = 51308.43  # H34: 51327.477 -> 51308.43
$ws.Cells.Item(34, 9).Value = 1183.8  # I34: 1263.8 -> 1183.8
$ws.Cells.Item(34, 11).Value = 1183.8  # K34: 1263.8 -> 1183.8
$ws.Cells.Item(34, 13).Value = -981.8  # M34: -1061.8 -> -981.8

# Row 58
$ws.Cells.Item(58, 8).Value = 4349.8335  # H58: 3856.7144 -> 4349.8335
$ws.Cells.Item(58, 9).Value = 5057  # I58: 4363.8335 -> 5057
$ws.Cells.Item(58, 11).Value = 5057  # K58: 4363.8335 -> 5057
$ws.Cells.Item(58, 13).Value = -4854  # M58: -4160.8335 -> -4854

# Row 122
$ws.Cells.Item(122, 8).Value = 4852.3335  # H122: 4248.4443 -> 4852.3335
$ws.Cells.Item(122, 9).Value = 6301  # I122: 5280.6 -> 6301
$ws.Cells.Item(122, 10).Value = 4325.5454  # J122: 3851.4614 -> 4325.5454
$ws.Cells.Item(122, 11).Value = 18903  # K122: 15841.8 -> 18903
$ws.Cells.Item(122, 12).Value = 12976.6362  # L122: 11554.3842 -> 12976.6362
$ws.Cells.Item(122, 13).Value = -16453  # M122: -13391.8 -> -16453
$ws.Cells.Item(122, 14).Value = -17876.6362  # N122: -16454.3842 -> -17876.6362

# Row 132
$ws.Cells.Item(132, 8).Value = 1491.6222  # H132: 1404.5834 -> 1491.6222
$ws.Cells.Item(132, 9).Value = 1223.4722  # I132: 1136.9744 -> 1223.4722
$ws.Cells.Item(132, 11).Value = 3670.4166  # K132: 3410.9232 -> 3670.4166
$ws.Cells.Item(132, 13).Value = -1140.4166  # M132: -880.9232000000002 -> -1140.4166

# Row 134
$ws.Cells.Item(134, 8).Value = 234689.94  # H134: 229379.7 -> 234689.94
$ws.Cells.Item(134, 9).Value = 2168.7073  # I134: 2136.1191 -> 2168.7073
$ws.Cells.Item(134, 10).Value = 5001375  # J134: 5001495 -> 5001375
$ws.Cells.Item(134, 11).Value = 6506.1219  # K134: 6408.3573 -> 6506.1219
$ws.Cells.Item(134, 12).Value = 15004125  # L134: 15004485 -> 15004125
$ws.Cells.Item(134, 13).Value = -3971.1219  # M134: -3873.3573 -> -3971.1219
$ws.Cells.Item(134, 14).Value = -15009195  # N134: -15009555 -> -15009195

# Row 136
$ws.Cells.Item(136, 8).Value = 4349.8335  # H136: 3856.7144 -> 4349.8335
$ws.Cells.Item(136, 9).Value = 5057  # I136: 4363.8335 -> 5057
$ws.Cells.Item(136, 11).Value = 15171  # K136: 13091.5005 -> 15171
$ws.Cells.Item(136, 13).Value = -12621  # M136: -10541.5005 -> -12621

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Cells.Item(131, 8).Value = 4554.76  # H131: 4219.3076 -> 4554.76
$ws.Cells.Item(131, 9).Value = 1613.6364  # I131: 1420 -> 1613.6364
$ws.Cells.Item(131, 10).Value = 6865.643  # J131: 7485.1665 -> 6865.643
$ws.Cells.Item(131, 11).Value = 4840.9092  # K131: 4260 -> 4840.9092
$ws.Cells.Item(131, 12).Value = 20596.929  # L131: 22455.4995 -> 20596.929
$ws.Cells.Item(131, 13).Value = 199.0907999999999  # M131: 780 -> 199.0907999999999
$ws.Cells.Item(131, 14).Value = -30676.929  # N131: -32535.4995 -> -30676.929

$ws = $wb.Worksheets.Item("GSM")
# Row 64
$ws.Cells.Item(64, 8).Value = 74771  # H64: 0 -> 74771
$ws.Cells.Item(64, 10).Value = 74771  # J64: 0 -> 74771
$ws.Cells.Item(64, 12).Value = 74771  # L64: 0 -> 74771
$ws.Cells.Item(64, 14).Value = -75267  # N64: None -> -75267

# Row 67
$ws.Cells.Item(67, 8).Value = 74771  # H67: 0 -> 74771
$ws.Cells.Item(67, 10).Value = 74771  # J67: 0 -> 74771
$ws.Cells.Item(67, 12).Value = 74771  # L67: 0 -> 74771
$ws.Cells.Item(67, 14).Value = -76487  # N67: None -> -76487

# Row 122
$ws.Cells.Item(122, 8).Value = 5479.9  # H122: 6099.875 -> 5479.9
$ws.Cells.Item(122, 9).Value = 3833.3333  # I122: 5500 -> 3833.3333
$ws.Cells.Item(122, 11).Value = 11499.9999  # K122: 16500 -> 11499.9999
$ws.Cells.Item(122, 13).Value = -9049.999899999999  # M122: -14050 -> -9049.999899999999

# Row 132
$ws.Cells.Item(132, 8).Value = 72350.94  # H132: 72350.19 -> 72350.94
$ws.Cells.Item(132, 9).Value = 9472.5  # I132: 9471.643 -> 9472.5
$ws.Cells.Item(132, 11).Value = 28417.5  # K132: 28414.929 -> 28417.5
$ws.Cells.Item(132, 13).Value = -25887.5  # M132: -25884.929 -> -25887.5

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Cells.Item(16, 8).Value = 0  # H16: 95.5 -> 0
$ws.Cells.Item(16, 9).Value = 0  # I16: 95.5 -> 0
$ws.Cells.Item(16, 11).Value = 0  # K16: 95.5 -> 0
$ws.Cells.Item(16, 13).Value = $null  # M16: clear (was 74.5)

# Row 100
$ws.Cells.Item(100, 8).Value = 3280.6  # H100: 2625.75 -> 3280.6
$ws.Cells.Item(100, 9).Value = 3099.75  # I100: 2625.75 -> 3099.75
$ws.Cells.Item(100, 10).Value = 4004  # J100: 0 -> 4004
$ws.Cells.Item(100, 11).Value = 3099.75  # K100: 2625.75 -> 3099.75
$ws.Cells.Item(100, 12).Value = 4004  # L100: 0 -> 4004
$ws.Cells.Item(100, 13).Value = -2558.75  # M100: -2084.75 -> -2558.75
$ws.Cells.Item(100, 14).Value = -5086  # N100: None -> -5086

# Row 122
$ws.Cells.Item(122, 8).Value = 5169.4614  # H122: 5491.5 -> 5169.4614
$ws.Cells.Item(122, 9).Value = 4356.222  # I122: 4700.2856 -> 4356.222
$ws.Cells.Item(122, 10).Value = 6999.25  # J122: 6599.2 -> 6999.25
$ws.Cells.Item(122, 11).Value = 13068.666  # K122: 14100.8568 -> 13068.666
$ws.Cells.Item(122, 12).Value = 20997.75  # L122: 19797.6 -> 20997.75
$ws.Cells.Item(122, 13).Value = -10618.666  # M122: -11650.8568 -> -10618.666
$ws.Cells.Item(122, 14).Value = -25897.75  # N122: -24697.6 -> -25897.75

# Row 132
$ws.Cells.Item(132, 8).Value = 7742.4546  # H132: 7666.8696 -> 7742.4546
$ws.Cells.Item(132, 9).Value = 6955.8667  # I132: 6896.375 -> 6955.8667
$ws.Cells.Item(132, 11).Value = 20867.6001  # K132: 20689.125 -> 20867.6001
$ws.Cells.Item(132, 13).Value = -18337.6001  # M132: -18159.125 -> -18337.6001

# Row 136
$ws.Cells.Item(136, 8).Value = 506210  # H136: 506240.66 -> 506210
$ws.Cells.Item(136, 9).Value = 1005299.1  # I136: 1116666.5 -> 1005299.1
$ws.Cells.Item(136, 10).Value = 7120.9  # J136: 6801.273 -> 7120.9
$ws.Cells.Item(136, 11).Value = 3015897.3  # K136: 3349999.5 -> 3015897.3
$ws.Cells.Item(136, 12).Value = 21362.7  # L136: 20403.819 -> 21362.7
$ws.Cells.Item(136, 13).Value = -3013347.3  # M136: -3347449.5 -> -3013347.3
$ws.Cells.Item(136, 14).Value = -26462.7  # N136: -25503.819 -> -26462.7

$ws = $wb.Worksheets.Item("WVR")
# Row 25
$ws.Cells.Item(25, 8).Value = 15000  # H25: 0 -> 15000
$ws.Cells.Item(25, 10).Value = 15000  # J25: 0 -> 15000
$ws.Cells.Item(25, 12).Value = 15000  # L25: 0 -> 15000
$ws.Cells.Item(25, 14).Value = -15586  # N25: None -> -15586

# Row 68
$ws.Cells.Item(68, 8).Value = 39943.5  # H68: 45000 -> 39943.5
$ws.Cells.Item(68, 10).Value = 39943.5  # J68: 45000 -> 39943.5
$ws.Cells.Item(68, 12).Value = 39943.5  # L68: 45000 -> 39943.5
$ws.Cells.Item(68, 14).Value = -41565.5  # N68: -46622 -> -41565.5

# Row 69
$ws.Cells.Item(69, 8).Value = 21954.2  # H69: 22208.4 -> 21954.2
$ws.Cells.Item(69, 10).Value = 21954.2  # J69: 22208.4 -> 21954.2
$ws.Cells.Item(69, 12).Value = 21954.2  # L69: 22208.4 -> 21954.2
$ws.Cells.Item(69, 14).Value = -23452.2  # N69: -23706.4 -> -23452.2

# Row 70
$ws.Cells.Item(70, 8).Value = 57517.5  # H70: 49385.625 -> 57517.5
$ws.Cells.Item(70, 9).Value = 25000  # I70: 24995 -> 25000
$ws.Cells.Item(70, 11).Value = 25000  # K70: 24995 -> 25000
$ws.Cells.Item(70, 13).Value = -24685  # M70: -24680 -> -24685

# Row 71
$ws.Cells.Item(71, 8).Value = 39943.5  # H71: 45000 -> 39943.5
$ws.Cells.Item(71, 10).Value = 39943.5  # J71: 45000 -> 39943.5
$ws.Cells.Item(71, 12).Value = 119830.5  # L71: 135000 -> 119830.5
$ws.Cells.Item(71, 14).Value = -127942.5  # N71: -143112 -> -127942.5

# Row 72
$ws.Cells.Item(72, 8).Value = 21954.2  # H72: 22208.4 -> 21954.2
$ws.Cells.Item(72, 10).Value = 21954.2  # J72: 22208.4 -> 21954.2
$ws.Cells.Item(72, 12).Value = 65862.60000000001  # L72: 66625.20000000001 -> 65862.60000000001
$ws.Cells.Item(72, 14).Value = -73350.60000000001  # N72: -74113.20000000001 -> -73350.60000000001

# Row 73
$ws.Cells.Item(73, 8).Value = 57517.5  # H73: 49385.625 -> 57517.5
$ws.Cells.Item(73, 9).Value = 25000  # I73: 24995 -> 25000
$ws.Cells.Item(73, 11).Value = 25000  # K73: 24995 -> 25000
$ws.Cells.Item(73, 13).Value = -23908  # M73: -23903 -> -23908

# Row 132
$ws.Cells.Item(132, 8).Value = 69118.31  # H132: 65505.47 -> 69118.31
$ws.Cells.Item(132, 10).Value = 107490.5  # J132: 98418.63 -> 107490.5
$ws.Cells.Item(132, 12).Value = 322471.5  # L132: 295255.89 -> 322471.5
$ws.Cells.Item(132, 14).Value = -327531.5  # N132: -300315.89 -> -327531.5

# Row 135
$ws.Cells.Item(135, 8).Value = 0  # H135: 50000 -> 0
$ws.Cells.Item(135, 10).Value = 0  # J135: 50000 -> 0
$ws.Cells.Item(135, 12).Value = 0  # L135: 50000 -> 0
$ws.Cells.Item(135, 14).Value = $null  # N135: clear (was -60140)

# Row 136
$ws.Cells.Item(136, 8).Value = 13975377  # H136: 13416394 -> 13975377
$ws.Cells.Item(136, 9).Value = 17547108  # I136: 16669793 -> 17547108
$ws.Cells.Item(136, 11).Value = 52641324  # K136: 50009379 -> 52641324
$ws.Cells.Item(136, 13).Value = -52638774  # M136: -50006829 -> -52638774

# Row 139
$ws.Cells.Item(139, 8).Value = 49000  # H139: 49666.668 -> 49000
$ws.Cells.Item(139, 10).Value = 0  # J139: 50000 -> 0
$ws.Cells.Item(139, 12).Value = 0  # L139: 50000 -> 0
$ws.Cells.Item(139, 14).Value = $null  # N139: clear (was -60280)

# Row 140
$ws.Cells.Item(140, 8).Value = 69000  # H140: 0 -> 69000
$ws.Cells.Item(140, 10).Value = 69000  # J140: 0 -> 69000
$ws.Cells.Item(140, 12).Value = 69000  # L140: 0 -> 69000
$ws.Cells.Item(140, 14).Value = -79360  # N140: None -> -79360

# Row 141
$ws.Cells.Item(141, 8).Value = 0  # H141: 50000 -> 0
$ws.Cells.Item(141, 10).Value = 0  # J141: 50000 -> 0
$ws.Cells.Item(141, 12).Value = 0  # L141: 50000 -> 0
$ws.Cells.Item(141, 14).Value = $null  # N141: clear (was -60360)
